{"js": "// Update the date paragraph and the 25 \"two-digit \u00f7 one-digit\" answer\n// cells inside the single table, per the target revision.\n//\n// Both the date text and every table-cell value are replaced by exact\n// *position* (paragraph index / row+column index) rather than by\n// searching for the old text, because one of the new values\n// (\"23\u00f78=2, 7\") happens to equal another cell's *old* value elsewhere\n// in the table - a plain global text replace could clobber it twice.\n\n// --- 1. Date paragraph (first paragraph in the body) ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\nif (dateParagraph.text.trim() === \"2025-12-05 Friday\") {\n  // Replace run text but keep paragraph formatting/run formatting by\n  // rewriting via insertText over the whole paragraph range.\n  dateParagraph.getRange().insertText(\"2025-12-06 Saturday\", Word.InsertLocation.replace);\n}\n\n// --- 2. Table of division answers ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Row-major list of old -> new values, exactly matching the five\n// populated rows (0, 4, 8, 12, 16) x five columns each.\nconst newGrid = [\n  [\"19\u00f79=2, 1\", \"96\u00f75=19, 1\"],\n  [\"58\u00f76=9, 4\", \"44\u00f77=6, 2\"],\n  [\"60\u00f72=30, 0\", \"62\u00f72=31, 0\"],\n  [\"77\u00f75=15, 2\", \"49\u00f73=16, 1\"],\n  [\"10\u00f74=2, 2\", \"79\u00f74=19, 3\"],\n  [\"32\u00f78=4, 0\", \"49\u00f74=12, 1\"],\n  [\"94\u00f72=47, 0\", \"57\u00f73=19, 0\"],\n  [\"96\u00f76=16, 0\", \"10\u00f78=1, 2\"],\n  [\"37\u00f73=12, 1\", \"83\u00f73=27, 2\"],\n  [\"82\u00f76=13, 4\", \"75\u00f78=9, 3\"],\n  [\"55\u00f72=27, 1\", \"30\u00f76=5, 0\"],\n  [\"77\u00f72=38, 1\", \"94\u00f74=23, 2\"],\n  [\"54\u00f76=9, 0\", \"91\u00f74=22, 3\"],\n  [\"76\u00f76=12, 4\", \"64\u00f72=32, 0\"],\n  [\"28\u00f79=3, 1\", \"70\u00f75=14, 0\"],\n  [\"97\u00f73=32, 1\", \"53\u00f77=7, 4\"],\n  [\"56\u00f78=7, 0\", \"22\u00f73=7, 1\"],\n  [\"81\u00f74=20, 1\", \"23\u00f78=2, 7\"],\n  [\"70\u00f79=7, 7\", \"45\u00f75=9, 0\"],\n  [\"13\u00f73=4, 1\", \"22\u00f73=7, 1\"],\n  [\"53\u00f72=26, 1\", \"35\u00f73=11, 2\"],\n  [\"54\u00f79=6, 0\", \"44\u00f79=4, 8\"],\n  [\"82\u00f72=41, 0\", \"82\u00f79=9, 1\"],\n  [\"75\u00f79=8, 3\", \"47\u00f76=7, 5\"],\n  [\"23\u00f78=2, 7\", \"55\u00f75=11, 0\"],\n];\n\nconst columnsPerRow = 5;\nlet idx = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < columnsPerRow; c++) {\n    const cell = table.getCellOrNullObject(r, c);\n    cell.load(\"value\");\n    await context.sync();\n    if (cell.isNullObject) continue;\n\n    const text = cell.value;\n    if (text === \"\") continue;\n\n    const pair = newGrid[idx];\n    if (pair && text === pair[0]) {\n      cell.value = pair[1];\n    }\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and the 25 \"two-digit \u00f7 one-digit\" answer\n# cells inside the single table, per the target revision.\n#\n# Both the date text and every table-cell value are addressed by exact\n# *position* (paragraph index / table row+column index) rather than by\n# searching for the old text, because one of the new values\n# (\"23\u00f78=2, 7\") happens to equal another cell's *old* value elsewhere\n# in the table - a plain global Find/Replace could clobber it twice.\n\n$d = $word.ActiveDocument\n\n# --- 1. Date paragraph (first paragraph in the body) ---\n$dateParagraph = $d.Paragraphs.Item(1)\nif ($dateParagraph.Range.Text.TrimEnd() -eq \"2025-12-05 Friday\") {\n    $dateParagraph.Range.Text = \"2025-12-06 Saturday\"\n}\n\n# --- 2. Table of division answers ---\n$tbl = $d.Tables.Item(1)\n\n# Row-major list of [old, new] values, exactly matching the five\n# populated rows (1,5,9,13,17) x five columns each, top-to-bottom,\n# left-to-right.\n$pairs = @(\n    @(\"19\u00f79=2, 1\", \"96\u00f75=19, 1\"),\n    @(\"58\u00f76=9, 4\", \"44\u00f77=6, 2\"),\n    @(\"60\u00f72=30, 0\", \"62\u00f72=31, 0\"),\n    @(\"77\u00f75=15, 2\", \"49\u00f73=16, 1\"),\n    @(\"10\u00f74=2, 2\", \"79\u00f74=19, 3\"),\n    @(\"32\u00f78=4, 0\", \"49\u00f74=12, 1\"),\n    @(\"94\u00f72=47, 0\", \"57\u00f73=19, 0\"),\n    @(\"96\u00f76=16, 0\", \"10\u00f78=1, 2\"),\n    @(\"37\u00f73=12, 1\", \"83\u00f73=27, 2\"),\n    @(\"82\u00f76=13, 4\", \"75\u00f78=9, 3\"),\n    @(\"55\u00f72=27, 1\", \"30\u00f76=5, 0\"),\n    @(\"77\u00f72=38, 1\", \"94\u00f74=23, 2\"),\n    @(\"54\u00f76=9, 0\", \"91\u00f74=22, 3\"),\n    @(\"76\u00f76=12, 4\", \"64\u00f72=32, 0\"),\n    @(\"28\u00f79=3, 1\", \"70\u00f75=14, 0\"),\n    @(\"97\u00f73=32, 1\", \"53\u00f77=7, 4\"),\n    @(\"56\u00f78=7, 0\", \"22\u00f73=7, 1\"),\n    @(\"81\u00f74=20, 1\", \"23\u00f78=2, 7\"),\n    @(\"70\u00f79=7, 7\", \"45\u00f75=9, 0\"),\n    @(\"13\u00f73=4, 1\", \"22\u00f73=7, 1\"),\n    @(\"53\u00f72=26, 1\", \"35\u00f73=11, 2\"),\n    @(\"54\u00f79=6, 0\", \"44\u00f79=4, 8\"),\n    @(\"82\u00f72=41, 0\", \"82\u00f79=9, 1\"),\n    @(\"75\u00f79=8, 3\", \"47\u00f76=7, 5\"),\n    @(\"23\u00f78=2, 7\", \"55\u00f75=11, 0\")\n)\n\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n$idx = 0\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cellText = $cell.Range.Text.TrimEnd(\"`r\", \"`a\")\n        if ($cellText -eq \"\") { continue }\n\n        $pair = $pairs[$idx]\n        if ($pair -and $cellText -eq $pair[0]) {\n            $cell.Range.Text = $pair[1]\n        }\n        $idx++\n    }\n}\n"}
